$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 6658
$ws.Range("F3").Value = 790
$ws.Range("F4").Value = 1101
$ws.Range("F8").Value = 7
$ws.Range("F9").Value = 1065
$ws.Range("F10").Value = 830
$ws.Range("F12").Value = 1309
$ws.Range("F16").Value = 535
$ws.Range("F17").Value = 7
$ws.Range("F18").Value = 362
$ws.Range("F20").Value = 1470
$ws.Range("F21").Value = 710
$ws.Range("F25").Value = 3
$ws.Range("F26").Value = 1119
$ws.Range("F27").Value = 240
$ws.Range("F28").Value = 2334
$ws.Range("F29").Value = 259
$ws.Range("F30").Value = 1182
$ws.Range("F33").Value = 3797

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 3
$ws.Range("F8").Value = 23
$ws.Range("F9").Value = 1031
$ws.Range("F11").Value = 147
$ws.Range("F18").Value = 329
$ws.Range("F19").Value = 4118
$ws.Range("F24").Value = 226
$ws.Range("F26").Value = 103
$ws.Range("G27").Value = 180
$ws.Range("F31").Value = 52
$ws.Range("F32").Value = 1698

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 1620
$ws.Range("F8").Value = 932

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1620
$ws.Range("F6").Value = 932
$ws.Range("F9").Value = 6658
$ws.Range("F10").Value = 3
$ws.Range("F12").Value = 790
$ws.Range("F16").Value = 1065
$ws.Range("F17").Value = 830
$ws.Range("F18").Value = 23
$ws.Range("F19").Value = 147
$ws.Range("F22").Value = 1309
$ws.Range("F24").Value = 524
$ws.Range("F25").Value = 535
$ws.Range("F26").Value = 329
$ws.Range("F27").Value = 362
$ws.Range("F28").Value = 1470
$ws.Range("F29").Value = 710
$ws.Range("F32").Value = 226
$ws.Range("F34").Value = 1119
$ws.Range("F35").Value = 240
$ws.Range("F36").Value = 103
$ws.Range("G37").Value = 180
$ws.Range("F38").Value = 2334
$ws.Range("F42").Value = 52
$ws.Range("F43").Value = 1698
$ws.Range("F44").Value = 1698
$ws.Range("F45").Value = 1182
$ws.Range("F48").Value = 3797
